$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Year"
$ws.Range("C1").Value = "Individual returns examined"
$ws.Range("D1").Value = "Individual exam coverage"
$ws.Range("G1").Value = "Individual percent of returns examined with no change (field)"
$ws.Range("H1").Value = "Individual percent of returns examined with no change (correspondence)"
$ws.Range("I1").Value = "Individual additional tax"
$ws.Range("J1").Value = "Individual additional tax (field)"
$ws.Range("K1").Value = "Individual additional tax (correspondence)"
$ws.Range("M1").Value = "Corporate returns examined"
$ws.Range("Q1").Value = "Corporate percent of returns examined with no change (field)"
$ws.Range("R1").Value = "Corporate percent of returns examined with no change (correspondence)"
$ws.Range("S1").Value = "Corporate additional tax"
$ws.Range("T1").Value = "Corporate additional tax (field)"
$ws.Range("U1").Value = "Corporate additional tax (correspondence)"

$ws.Range("F16").Select()
